$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.146379828453064
$ws.Range("B1").Value = 2.570617198944092
$ws.Range("C1").Value = 2.889833927154541
$ws.Range("D1").Value = 2.895297050476074
$ws.Range("E1").Value = 0.6916077136993408
